$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("figures")

# The table gained a new "r-generated" indicator column (it previously showed
# up as an auto-named "Column1" with only a sparse free-text note in E7).
# Rename the header and populate the indicator for every figure/table row:
# rows 2-7 (not produced with R) get 0, the rest (R-generated) get 1.
$ws.Range("E1").Value = "r-generated"

for ($row = 2; $row -le 30; $row++) {
    if ($row -le 7) {
        $ws.Cells.Item($row, 5).Value = 0
    } else {
        $ws.Cells.Item($row, 5).Value = 1
    }
}

# Set up the page for printing this sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Reflect the author's last working cell selection.
$ws.Range("C10").Select()

$wb.Save()
